$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style (bold, border, centered) from the last existing index cell (A155)
# into the new index-column cells for rows 156-161, matching the sheet's
# existing formatting convention for column A.
$ws.Range("A155").Copy() | Out-Null
$ws.Range("A156:A161").PasteSpecial(-4122) | Out-Null

# Populate the six new box-score rows (154-159 in the "index" column A,
# spreadsheet rows 156-161) with the latest game data.
$ws.Cells.Item(156, 1).Value = 154
$ws.Cells.Item(156, 2).Value = "HOU"
$ws.Cells.Item(156, 3).Value = "GSW"
$ws.Cells.Item(156, 4).Value = "away"
$ws.Cells.Item(156, 5).NumberFormat = "@"
$ws.Cells.Item(156, 5).Value = "2025-04-06"
$ws.Cells.Item(156, 5).ClearFormats()
$ws.Cells.Item(156, 6).Value = "240:00"
$ws.Cells.Item(156, 7).Value = 45
$ws.Cells.Item(156, 8).Value = 93
$ws.Cells.Item(156, 9).Value = 0.484
$ws.Cells.Item(156, 10).Value = 10
$ws.Cells.Item(156, 11).Value = 30
$ws.Cells.Item(156, 12).Value = 0.333
$ws.Cells.Item(156, 13).Value = 6
$ws.Cells.Item(156, 14).Value = 16
$ws.Cells.Item(156, 15).Value = 0.375
$ws.Cells.Item(156, 16).Value = 16
$ws.Cells.Item(156, 17).Value = 33
$ws.Cells.Item(156, 18).Value = 49
$ws.Cells.Item(156, 19).Value = 26
$ws.Cells.Item(156, 20).Value = 11
$ws.Cells.Item(156, 21).Value = 7
$ws.Cells.Item(156, 22).Value = 15
$ws.Cells.Item(156, 23).Value = 15
$ws.Cells.Item(156, 24).Value = 106
$ws.Cells.Item(156, 25).Value = 10
$ws.Cells.Item(156, 26).Value = 18
$ws.Cells.Item(156, 27).Value = 33
$ws.Cells.Item(156, 28).Value = 32
$ws.Cells.Item(156, 29).Value = 23
$ws.Cells.Item(156, 30).Value = "W"
$ws.Cells.Item(157, 1).Value = 155
$ws.Cells.Item(157, 2).Value = "GSW"
$ws.Cells.Item(157, 3).Value = "HOU"
$ws.Cells.Item(157, 4).Value = "home"
$ws.Cells.Item(157, 5).NumberFormat = "@"
$ws.Cells.Item(157, 5).Value = "2025-04-06"
$ws.Cells.Item(157, 5).ClearFormats()
$ws.Cells.Item(157, 6).Value = "240:00"
$ws.Cells.Item(157, 7).Value = 37
$ws.Cells.Item(157, 8).Value = 89
$ws.Cells.Item(157, 9).Value = 0.416
$ws.Cells.Item(157, 10).Value = 16
$ws.Cells.Item(157, 11).Value = 43
$ws.Cells.Item(157, 12).Value = 0.372
$ws.Cells.Item(157, 13).Value = 6
$ws.Cells.Item(157, 14).Value = 11
$ws.Cells.Item(157, 15).Value = 0.545
$ws.Cells.Item(157, 16).Value = 15
$ws.Cells.Item(157, 17).Value = 31
$ws.Cells.Item(157, 18).Value = 46
$ws.Cells.Item(157, 19).Value = 28
$ws.Cells.Item(157, 20).Value = 12
$ws.Cells.Item(157, 21).Value = 2
$ws.Cells.Item(157, 22).Value = 20
$ws.Cells.Item(157, 23).Value = 14
$ws.Cells.Item(157, 24).Value = 96
$ws.Cells.Item(157, 25).Value = -10
$ws.Cells.Item(157, 26).Value = 24
$ws.Cells.Item(157, 27).Value = 26
$ws.Cells.Item(157, 28).Value = 26
$ws.Cells.Item(157, 29).Value = 20
$ws.Cells.Item(157, 30).Value = "L"
$ws.Cells.Item(158, 1).Value = 156
$ws.Cells.Item(158, 2).Value = "GSW"
$ws.Cells.Item(158, 3).Value = "PHX"
$ws.Cells.Item(158, 4).Value = "away"
$ws.Cells.Item(158, 5).NumberFormat = "@"
$ws.Cells.Item(158, 5).Value = "2025-04-08"
$ws.Cells.Item(158, 5).ClearFormats()
$ws.Cells.Item(158, 6).Value = "240:00"
$ws.Cells.Item(158, 7).Value = 45
$ws.Cells.Item(158, 8).Value = 93
$ws.Cells.Item(158, 9).Value = 0.484
$ws.Cells.Item(158, 10).Value = 15
$ws.Cells.Item(158, 11).Value = 40
$ws.Cells.Item(158, 12).Value = 0.375
$ws.Cells.Item(158, 13).Value = 28
$ws.Cells.Item(158, 14).Value = 34
$ws.Cells.Item(158, 15).Value = 0.824
$ws.Cells.Item(158, 16).Value = 19
$ws.Cells.Item(158, 17).Value = 38
$ws.Cells.Item(158, 18).Value = 57
$ws.Cells.Item(158, 19).Value = 31
$ws.Cells.Item(158, 20).Value = 14
$ws.Cells.Item(158, 21).Value = 3
$ws.Cells.Item(158, 22).Value = 15
$ws.Cells.Item(158, 23).Value = 10
$ws.Cells.Item(158, 24).Value = 133
$ws.Cells.Item(158, 25).Value = 38
$ws.Cells.Item(158, 26).Value = 37
$ws.Cells.Item(158, 27).Value = 32
$ws.Cells.Item(158, 28).Value = 26
$ws.Cells.Item(158, 29).Value = 38
$ws.Cells.Item(158, 30).Value = "W"
$ws.Cells.Item(159, 1).Value = 157
$ws.Cells.Item(159, 2).Value = "PHX"
$ws.Cells.Item(159, 3).Value = "GSW"
$ws.Cells.Item(159, 4).Value = "home"
$ws.Cells.Item(159, 5).NumberFormat = "@"
$ws.Cells.Item(159, 5).Value = "2025-04-08"
$ws.Cells.Item(159, 5).ClearFormats()
$ws.Cells.Item(159, 6).Value = "240:00"
$ws.Cells.Item(159, 7).Value = 35
$ws.Cells.Item(159, 8).Value = 88
$ws.Cells.Item(159, 9).Value = 0.398
$ws.Cells.Item(159, 10).Value = 12
$ws.Cells.Item(159, 11).Value = 43
$ws.Cells.Item(159, 12).Value = 0.279
$ws.Cells.Item(159, 13).Value = 13
$ws.Cells.Item(159, 14).Value = 15
$ws.Cells.Item(159, 15).Value = 0.867
$ws.Cells.Item(159, 16).Value = 13
$ws.Cells.Item(159, 17).Value = 28
$ws.Cells.Item(159, 18).Value = 41
$ws.Cells.Item(159, 19).Value = 22
$ws.Cells.Item(159, 20).Value = 12
$ws.Cells.Item(159, 21).Value = 4
$ws.Cells.Item(159, 22).Value = 19
$ws.Cells.Item(159, 23).Value = 19
$ws.Cells.Item(159, 24).Value = 95
$ws.Cells.Item(159, 25).Value = -38
$ws.Cells.Item(159, 26).Value = 24
$ws.Cells.Item(159, 27).Value = 19
$ws.Cells.Item(159, 28).Value = 18
$ws.Cells.Item(159, 29).Value = 34
$ws.Cells.Item(159, 30).Value = "L"
$ws.Cells.Item(160, 1).Value = 158
$ws.Cells.Item(160, 2).Value = "SAS"
$ws.Cells.Item(160, 3).Value = "GSW"
$ws.Cells.Item(160, 4).Value = "away"
$ws.Cells.Item(160, 5).NumberFormat = "@"
$ws.Cells.Item(160, 5).Value = "2025-04-09"
$ws.Cells.Item(160, 5).ClearFormats()
$ws.Cells.Item(160, 6).Value = "240:00"
$ws.Cells.Item(160, 7).Value = 39
$ws.Cells.Item(160, 8).Value = 81
$ws.Cells.Item(160, 9).Value = 0.481
$ws.Cells.Item(160, 10).Value = 18
$ws.Cells.Item(160, 11).Value = 46
$ws.Cells.Item(160, 12).Value = 0.391
$ws.Cells.Item(160, 13).Value = 18
$ws.Cells.Item(160, 14).Value = 21
$ws.Cells.Item(160, 15).Value = 0.857
$ws.Cells.Item(160, 16).Value = 10
$ws.Cells.Item(160, 17).Value = 33
$ws.Cells.Item(160, 18).Value = 43
$ws.Cells.Item(160, 19).Value = 23
$ws.Cells.Item(160, 20).Value = 6
$ws.Cells.Item(160, 21).Value = 5
$ws.Cells.Item(160, 22).Value = 13
$ws.Cells.Item(160, 23).Value = 23
$ws.Cells.Item(160, 24).Value = 114
$ws.Cells.Item(160, 25).Value = 3
$ws.Cells.Item(160, 26).Value = 23
$ws.Cells.Item(160, 27).Value = 32
$ws.Cells.Item(160, 28).Value = 21
$ws.Cells.Item(160, 29).Value = 38
$ws.Cells.Item(160, 30).Value = "W"
$ws.Cells.Item(161, 1).Value = 159
$ws.Cells.Item(161, 2).Value = "GSW"
$ws.Cells.Item(161, 3).Value = "SAS"
$ws.Cells.Item(161, 4).Value = "home"
$ws.Cells.Item(161, 5).NumberFormat = "@"
$ws.Cells.Item(161, 5).Value = "2025-04-09"
$ws.Cells.Item(161, 5).ClearFormats()
$ws.Cells.Item(161, 6).Value = "240:00"
$ws.Cells.Item(161, 7).Value = 36
$ws.Cells.Item(161, 8).Value = 86
$ws.Cells.Item(161, 9).Value = 0.419
$ws.Cells.Item(161, 10).Value = 16
$ws.Cells.Item(161, 11).Value = 48
$ws.Cells.Item(161, 12).Value = 0.333
$ws.Cells.Item(161, 13).Value = 23
$ws.Cells.Item(161, 14).Value = 28
$ws.Cells.Item(161, 15).Value = 0.821
$ws.Cells.Item(161, 16).Value = 16
$ws.Cells.Item(161, 17).Value = 29
$ws.Cells.Item(161, 18).Value = 45
$ws.Cells.Item(161, 19).Value = 29
$ws.Cells.Item(161, 20).Value = 8
$ws.Cells.Item(161, 21).Value = 2
$ws.Cells.Item(161, 22).Value = 10
$ws.Cells.Item(161, 23).Value = 19
$ws.Cells.Item(161, 24).Value = 111
$ws.Cells.Item(161, 25).Value = -3
$ws.Cells.Item(161, 26).Value = 32
$ws.Cells.Item(161, 27).Value = 19
$ws.Cells.Item(161, 28).Value = 37
$ws.Cells.Item(161, 29).Value = 23
$ws.Cells.Item(161, 30).Value = "L"
